$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(1)

# Fill in row 3 (autogluon) values
$ws.Range("B3").Value = "0.339 (0.292 ± 0.020)"
$ws.Range("C3").Value = "00:03:03 (00:03:31 ± 00:00:25)"
$ws.Range("D3").Value = "00:00:00 (00:00:00 ± 00:00:00)"
$ws.Range("E3").Value = "[]"

# F3 should be stored as text "43" (matching sibling cells F4/F6/F8),
# not a numeric value, so force text formatting before/after assignment.
$ws.Range("F3").NumberFormat = "@"
$ws.Range("F3").Value = "43"
$ws.Range("F3").Style = "Normal"

# Fix mojibake "Â±" -> "±" in rows 4, 6, 8 (columns B, C, D)
$ws.Range("B4").Value = "0.758 (0.728 ± 0.020)"
$ws.Range("C4").Value = "00:01:07 (00:01:27 ± 00:00:19)"
$ws.Range("D4").Value = "00:00:01 (00:00:01 ± 00:00:00)"

$ws.Range("B6").Value = "0.816 (0.791 ± 0.016)"
$ws.Range("C6").Value = "00:04:56 (00:05:00 ± 00:00:02)"
$ws.Range("D6").Value = "00:00:00 (00:00:04 ± 00:00:02)"

$ws.Range("B8").Value = "0.781 (0.667 ± 0.037)"
$ws.Range("C8").Value = "00:04:51 (00:07:35 ± 00:02:34)"
$ws.Range("D8").Value = "00:00:00 (00:00:00 ± 00:00:00)"
